$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.389.80"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.943.25"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.60"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.03"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9754"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "1.935.58"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.091"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.761"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07055"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.75"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009825"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.09"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "29.413.83"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.479"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "2.171.34"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.98"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.752"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.45"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09343"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8613"
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.176"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.307"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.105"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05775"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.154"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.699"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5678"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1781"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.411"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.735"
$ws.Range("E44").Value = "  +7.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002820"
$ws.Range("E45").Value = "  +33.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5295"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.43"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06869"
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.081"
$ws.Range("E49").Value = "  -5.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.817"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  -1.59%  "
